# Commit: "Cleaning the Networks for analysis"
# Fix hyphenated words ("-level", "-century", "Union-27") in the reference
# titles so they read with a normal space instead of a hyphen.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("A20").Value = "Vicente M 2011-1 - Assessing the regional digital divide across the European Union 27"
$ws.Range("A11").Value = "Van Deursen A 2017-1 The relation between 21st century skills and digital skills: A systematic literature review"
$ws.Range("A10").Value = "Van Deursen A 2019 - The first level digital divide shifts from inequalities in physical access to inequalities in material access"
$ws.Range("A7").Value  = "Scheerder A 2017 -Determinants of Internet skills, uses and outcomes. A systematic review of the second  and third level digital divide"
$ws.Range("A6").Value  = "Hargittai E 2002 - Second Level Digital Divide: Differences in People's Online Skills"

# Restore the user's final selection to the data range, as in the saved file.
$ws.Range("A2:E32").Select()
